# Trade #12 closed at 2026-02-17 19:47:49 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" sheet metrics, the "Strategy Status" row for
# MarketMaking, and appends the new closed trade as row 13 on both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1300.27
$summary.Range("B4").Value = 0.27
$summary.Range("B5").Value = 0.45
$summary.Range("B6").Value = 12
$summary.Range("B7").Value = 7
$summary.Range("B9").Value = 58.33

# --- Strategy Status sheet (MarketMaking row) -----------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.27
$status.Range("D4").Value = 12
$status.Range("E4").Value = 0.27
$status.Range("F4").Value = 0.27
$status.Range("G4").Value = 58.33

# --- Append new trade row (row 13) to "All Trades" and "MarketMaking" ----
# Start from a copy of the last existing row (row 12) so that text-like
# values (dates, times, strategy/side/status labels) keep their original
# string typing/formatting instead of being auto-converted by Excel's
# smart-entry parsing (e.g. assigning "2026-02-17" to .Value turns it into
# a date serial number), then overwrite only the cells whose values
# actually differ from row 12. Column B (Date) is already "2026-02-17"
# after the copy, same as the new trade's date, so it is left untouched.
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A12:Q12").Copy()
    $ws.Range("A13").PasteSpecial()

    $ws.Cells.Item(13, 1).Value = 12              # Trade #
    $ws.Cells.Item(13, 3).Value = "19:47:43"       # Time
    $ws.Cells.Item(13, 6).Value = 0.44             # Entry Price
    $ws.Cells.Item(13, 7).Value = 0.48             # Exit Price
    $ws.Cells.Item(13, 9).Value = 9.0909           # P&L %
    $ws.Cells.Item(13, 10).Value = 0.04            # P&L $
    $ws.Cells.Item(13, 11).Value = 100.27          # Capital After
    $ws.Cells.Item(13, 17).Value = 0.13
}
